$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Monsters")

# Row 332: keep AF:AK as-is, just add AL:AO equal to AK's value
$val332 = $ws.Range("AK332").Value2
$ws.Range("AL332").Value2 = $val332
$ws.Range("AM332").Value2 = $val332
$ws.Range("AN332").Value2 = $val332
$ws.Range("AO332").Value2 = $val332

# Rows 333 - 375: replace AF:AJ with AK's original value, then extend AK's
# value across AL:AO as well.
for ($r = 333; $r -le 375; $r++) {
    $akVal = $ws.Range("AK$r").Value2
    $ws.Range("AF$r").Value2 = $akVal
    $ws.Range("AG$r").Value2 = $akVal
    $ws.Range("AH$r").Value2 = $akVal
    $ws.Range("AI$r").Value2 = $akVal
    $ws.Range("AJ$r").Value2 = $akVal
    $ws.Range("AL$r").Value2 = $akVal
    $ws.Range("AM$r").Value2 = $akVal
    $ws.Range("AN$r").Value2 = $akVal
    $ws.Range("AO$r").Value2 = $akVal
}

# Update the sheet view scroll position / selection to match the edit
$ws.Application.ActiveWindow.ScrollRow = 319
$ws.Application.ActiveWindow.ScrollColumn = 33
$ws.Range("AN323").Select()
